$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 'Add position'
$ws.Range("B3").Value = 'FAILED'
$ws.Range("A4").Value = 'Add position'
$ws.Range("B4").Value = 'FAILED'
$ws.Range("A5").Value = 'Add position'
$ws.Range("B5").Value = 'FAILED'
$ws.Range("A6").Value = 'Add position'
$ws.Range("B6").Value = 'FAILED'
$ws.Range("A7").Value = 'Add position'
$ws.Range("B7").Value = 'FAILED'
$ws.Range("A8").Value = 'Add position'
$ws.Range("B8").Value = 'FAILED'
$ws.Range("A9").Value = 'Add position'
$ws.Range("B9").Value = 'FAILED'
$ws.Range("A10").Value = 'Add position'
$ws.Range("B10").Value = 'FAILED'
$ws.Range("A11").Value = 'Add position'
$ws.Range("B11").Value = 'FAILED'
$ws.Range("A12").Value = 'Add position'
$ws.Range("B12").Value = 'FAILED'
$ws.Range("A13").Value = 'Add position'
$ws.Range("B13").Value = 'FAILED'
$ws.Range("A14").Value = 'Add position'
$ws.Range("B14").Value = 'FAILED'
$ws.Range("A15").Value = 'Add position'
$ws.Range("B15").Value = 'FAILED'
$ws.Range("A16").Value = 'Add position'
$ws.Range("B16").Value = 'FAILED'
$ws.Range("A17").Value = 'Add position in Document Types'
$ws.Range("B17").Value = 'FAILED'
$ws.Range("A18").Value = 'Add position in Document Types'
$ws.Range("B18").Value = 'FAILED'
$ws.Range("A19").Value = 'Add position in Document Types'
$ws.Range("B19").Value = 'FAILED'
$ws.Range("A20").Value = 'Add position in Document Types'
$ws.Range("B20").Value = 'FAILED'
$ws.Range("A21").Value = 'Add position in Document Types'
$ws.Range("B21").Value = 'FAILED'
$ws.Range("A22").Value = 'Add position in Document Types'
$ws.Range("B22").Value = 'FAILED'
$ws.Range("A23").Value = 'Add position in Document Types'
$ws.Range("B23").Value = 'FAILED'
$ws.Range("A24").Value = 'Add position in Document Types'
$ws.Range("B24").Value = 'FAILED'
$ws.Range("A25").Value = 'Add position in Document Types'
$ws.Range("B25").Value = 'FAILED'
$ws.Range("A26").Value = 'Add position in Document Types'
$ws.Range("B26").Value = 'FAILED'
$ws.Range("A27").Value = 'Add position in Document Types'
$ws.Range("B27").Value = 'FAILED'
$ws.Range("A28").Value = 'Add position in Document Types'
$ws.Range("B28").Value = 'FAILED'
$ws.Range("A29").Value = 'Add position in Document Types'
$ws.Range("B29").Value = 'FAILED'
$ws.Range("A30").Value = 'Add position in Document Types'
$ws.Range("B30").Value = 'FAILED'
$ws.Range("A31").Value = 'Add position in Document Types'
$ws.Range("B31").Value = 'FAILED'
$ws.Range("A32").Value = 'Add position in Document Types'
$ws.Range("B32").Value = 'FAILED'
$ws.Range("A33").Value = 'User adds, edits and deletes in Fields'
$ws.Range("B33").Value = 'FAILED'
$ws.Range("A34").Value = 'User adds, edits and deletes in Fields'
$ws.Range("B34").Value = 'FAILED'
$ws.Range("A35").Value = 'User adds, edits and deletes in Fields'
$ws.Range("B35").Value = 'FAILED'
$ws.Range("A36").Value = 'Add position in Fields'
$ws.Range("B36").Value = 'FAILED'
$ws.Range("A37").Value = 'Add position in Fields'
$ws.Range("B37").Value = 'FAILED'
$ws.Range("A38").Value = 'Add position in Fields'
$ws.Range("B38").Value = 'FAILED'
$ws.Range("A39").Value = 'Add position in Fields'
$ws.Range("B39").Value = 'FAILED'
$ws.Range("A40").Value = 'Add position in Fields'
$ws.Range("B40").Value = 'FAILED'
$ws.Range("A41").Value = 'Add position in Fields'
$ws.Range("B41").Value = 'FAILED'
$ws.Range("A42").Value = 'User adds, edits and deletes data in the Position Management functionality'
$ws.Range("B42").Value = 'PASSED'
$ws.Range("A43").Value = 'Add position in Fields'
$ws.Range("B43").Value = 'FAILED'
$ws.Range("A44").Value = 'Add position in Fields'
$ws.Range("B44").Value = 'FAILED'
$ws.Range("A45").Value = 'Add position in Fields'
$ws.Range("B45").Value = 'FAILED'
$ws.Range("A46").Value = 'Add position in Fields'
$ws.Range("B46").Value = 'FAILED'
$ws.Range("A47").Value = 'Add position in Fields'
$ws.Range("B47").Value = 'FAILED'
$ws.Range("A48").Value = 'Add position in Fields'
$ws.Range("B48").Value = 'FAILED'
$ws.Range("A49").Value = 'Add position in Fields'
$ws.Range("B49").Value = 'FAILED'
$ws.Range("A50").Value = 'Add position in Fields'
$ws.Range("B50").Value = 'FAILED'
$ws.Range("A51").Value = 'Add position in Fields'
$ws.Range("B51").Value = 'FAILED'
$ws.Range("A52").Value = 'Add position in Fields'
$ws.Range("B52").Value = 'FAILED'
$ws.Range("A53").Value = 'Add position in Document Types'
$ws.Range("B53").Value = 'FAILED'
$ws.Range("A54").Value = 'Add position in Document Types'
$ws.Range("B54").Value = 'FAILED'
$ws.Range("A55").Value = 'Add position in Document Types'
$ws.Range("B55").Value = 'FAILED'
$ws.Range("A56").Value = 'Add position in Document Types'
$ws.Range("B56").Value = 'FAILED'
$ws.Range("A57").Value = 'Add position in Document Types'
$ws.Range("B57").Value = 'FAILED'
$ws.Range("A58").Value = 'Add position in Document Types'
$ws.Range("B58").Value = 'FAILED'
$ws.Range("A59").Value = 'Add position in Document Types'
$ws.Range("B59").Value = 'FAILED'
$ws.Range("A60").Value = 'Add position in Document Types'
$ws.Range("B60").Value = 'PASSED'
$ws.Range("A61").Value = 'Add position in Document Types'
$ws.Range("B61").Value = 'FAILED'
$ws.Range("A62").Value = 'Add position in Document Types'
$ws.Range("B62").Value = 'FAILED'
$ws.Range("A63").Value = 'User adds, edits and deletes data in the Attestations functionality'
$ws.Range("B63").Value = 'FAILED'
$ws.Range("A64").Value = 'Add position in Document Types'
$ws.Range("B64").Value = 'FAILED'